$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (rows 3-21): change the date separators from "/" to "-" ---
# Force the range to Text format first so Excel does not reinterpret the
# dash-separated string as a date value, then restore the default cell
# style afterwards so no stray number formatting is left on the cells.
$dateRange = $ws.Range("A3:A21")
$dateRange.NumberFormat = "@"

for ($r = 3; $r -le 21; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $oldVal = $cell.Value()
    $newVal = $oldVal.Replace("/", "-")
    $cell.Value = $newVal
}

$dateRange.Style = "Normal"

# --- Update the attendance flag cells (columns D, E, G, H) ---
$ws.Cells.Item(3, 4).Value = 1   # D3
$ws.Cells.Item(3, 7).Value = 1   # G3

$ws.Cells.Item(4, 4).Value = 1   # D4
$ws.Cells.Item(4, 5).Value = 1   # E4
$ws.Cells.Item(4, 8).Value = 0   # H4

$ws.Cells.Item(5, 4).Value = 1   # D5
$ws.Cells.Item(5, 5).Value = 1   # E5
$ws.Cells.Item(5, 8).Value = 0   # H5

$ws.Cells.Item(10, 4).Value = 1  # D10
$ws.Cells.Item(10, 5).Value = 1  # E10
$ws.Cells.Item(10, 8).Value = 0  # H10

$ws.Cells.Item(14, 4).Value = 1  # D14
$ws.Cells.Item(14, 5).Value = 1  # E14
$ws.Cells.Item(14, 8).Value = 0  # H14
